# Updates the crypto price/volume table on Sheet1 with the latest scrape
# values (GitHub Actions refresh). Price column (D) holds plain text that
# looks numeric (e.g. "25.792.81", "0.0₃0776") so each assignment is
# apostrophe-prefixed to force Excel to keep it as text instead of
# auto-converting/parsing it as a number (which would corrupt values like
# "0.550" -> 0.55 or multi-dot prices). The Volume(1h) column (E) already
# contains padded, percent-suffixed text so it round-trips as text as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.792.81"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").Value = "'1.636.36"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'215.24"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("E6").Value = "  -0.76%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").Value = "'19.86"
$ws.Range("E10").Value = "  +1.39%  "

$ws.Range("D11").Value = "'0.0785"
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("E12").Value = "  -0.85%  "

$ws.Range("D13").Value = "'1.642.44"
$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("D14").Value = "'1.861.77"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").Value = "'0.0₃0776"
$ws.Range("E16").Value = "  +1.88%  "

$ws.Range("D17").Value = "'63.05"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "'25.817.59"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("E20").Value = "  +2.70%  "

$ws.Range("D21").Value = "'193.85"
$ws.Range("E21").Value = "  -0.93%  "

$ws.Range("D22").Value = "'9.94"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").Value = "'6.16"
$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "'139.25"
$ws.Range("E26").Value = "  -0.47%  "

$ws.Range("E27").Value = "  -4.79%  "

$ws.Range("D28").Value = "'6.83"
$ws.Range("E28").Value = "  +0.52%  "

$ws.Range("D29").Value = "'15.55"
$ws.Range("E29").Value = "  +0.50%  "

$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("D31").Value = "'0.0496"
$ws.Range("E31").Value = "  +1.56%  "

$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("E34").Value = "  +2.52%  "

$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("E36").Value = "  -0.66%  "

$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "'0.550"
$ws.Range("E38").Value = "  -0.24%  "

$ws.Range("D39").Value = "'1.108.96"
$ws.Range("E39").Value = "  -1.90%  "

$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("E41").Value = "  +0.76%  "

$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("D43").Value = "'99.21"
$ws.Range("E43").Value = "  +1.51%  "

$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("D45").Value = "'0.0₆0109"
$ws.Range("E45").Value = "  -3.20%  "

$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").Value = "'2.55"
$ws.Range("E46").Value = "  +14.01%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.61"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'7.73"
$ws.Range("E48").Value = "  +0.26%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.418"
$ws.Range("E49").Value = "  -5.99%  "

$ws.Range("E50").Value = "  -0.37%  "

$ws.Range("E51").Value = "  +0.07%  "
